$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, even if it looks numeric
# (e.g. "329.04"), without leaving a lasting number-format override
# on the cell once done.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '27.934.83'
$ws.Range("E2").Value = '  +0.72%  '

$ws.Range("D3").Value = '1.767.96'
$ws.Range("E3").Value = '  -0.31%  '

Set-TextValue $ws.Range("D4") '1.002'
$ws.Range("E4").Value = '  -0.10%  '

Set-TextValue $ws.Range("D5") '329.04'
$ws.Range("E5").Value = '  +0.60%  '

Set-TextValue $ws.Range("D6") '1.001'
$ws.Range("E6").Value = '  -0.14%  '

Set-TextValue $ws.Range("D7") '0.4683'
$ws.Range("E7").Value = '  +2.29%  '

Set-TextValue $ws.Range("D8") '0.3524'

Set-TextValue $ws.Range("D9") '43.66'
$ws.Range("E9").Value = '  +4.61%  '

Set-TextValue $ws.Range("D10") '0.07383'
$ws.Range("E10").Value = '  -1.28%  '

Set-TextValue $ws.Range("D11") '1.084'
$ws.Range("E11").Value = '  -1.57%  '

Set-TextValue $ws.Range("D12") '1.002'
$ws.Range("E12").Value = '  -0.06%  '

Set-TextValue $ws.Range("D13") '20.63'
$ws.Range("E13").Value = '  -0.86%  '

Set-TextValue $ws.Range("D14") '6.016'
$ws.Range("E14").Value = '  -0.25%  '

Set-TextValue $ws.Range("D15") '7.194'
$ws.Range("E15").Value = '  -0.10%  '

$ws.Range("D16").Value = '1.767.17'
$ws.Range("E16").Value = '  -0.31%  '

Set-TextValue $ws.Range("D17") '92.19'
$ws.Range("E17").Value = '  -1.38%  '

Set-TextValue $ws.Range("D18") '0.00001056'
$ws.Range("E18").Value = '  -0.20%  '

Set-TextValue $ws.Range("D19") '0.06417'
$ws.Range("E19").Value = '  -0.22%  '

Set-TextValue $ws.Range("D21") '16.92'
$ws.Range("E21").Value = '  -0.70%  '

Set-TextValue $ws.Range("D22") '5.799'
$ws.Range("E22").Value = '  +0.06%  '

$ws.Range("D23").Value = '27.989.94'
$ws.Range("E23").Value = '  +0.77%  '

Set-TextValue $ws.Range("D24") '11.15'
$ws.Range("E24").Value = '  -1.36%  '

Set-TextValue $ws.Range("D25") '2.157'
$ws.Range("E25").Value = '  +3.51%  '

Set-TextValue $ws.Range("D26") '165.09'
$ws.Range("E26").Value = '  +0.61%  '

$ws.Range("E27").Value = '  -0.72%  '

$ws.Range("D28").Value = '1.970.23'
$ws.Range("E28").Value = '  -0.37%  '

Set-TextValue $ws.Range("D29") '2.206'
$ws.Range("E29").Value = '  +1.53%  '

Set-TextValue $ws.Range("D30") '123.40'
$ws.Range("E30").Value = '  -1.68%  '

Set-TextValue $ws.Range("D31") '1.075'
$ws.Range("E31").Value = '  -2.16%  '

Set-TextValue $ws.Range("D32") '0.09341'
$ws.Range("E32").Value = '  +1.55%  '

Set-TextValue $ws.Range("D33") '3.656'
$ws.Range("E33").Value = '  -0.43%  '

Set-TextValue $ws.Range("D34") '5.552'
$ws.Range("E34").Value = '  +0.45%  '

Set-TextValue $ws.Range("D35") '11.68'
$ws.Range("E35").Value = '  -1.33%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D36") '0.06107'
$ws.Range("E36").Value = '  -0.93%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D37") '0.02263'
$ws.Range("E37").Value = '  -1.19%  '

Set-TextValue $ws.Range("D38") '0.2071'
$ws.Range("E38").Value = '  -0.87%  '

Set-TextValue $ws.Range("D39") '4.907'
$ws.Range("E39").Value = '  -0.87%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D40") '0.6178'
$ws.Range("E40").Value = '  -2.18%  '

$ws.Range("B41").Value = 'WEMIXTOKEN'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D41") '1.447'
$ws.Range("E41").Value = '  +4.11%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D42") '1.192'
$ws.Range("E42").Value = '  +0.65%  '

Set-TextValue $ws.Range("D43") '7.758'
$ws.Range("E43").Value = '  -0.62%  '

Set-TextValue $ws.Range("D44") '13.16'
$ws.Range("E44").Value = '  -0.78%  '

Set-TextValue $ws.Range("D45") '3.750'
$ws.Range("E45").Value = '  +0.36%  '

Set-TextValue $ws.Range("D46") '0.5801'
$ws.Range("E46").Value = '  -1.90%  '

Set-TextValue $ws.Range("D47") '124.02'
$ws.Range("E47").Value = '  +1.12%  '

Set-TextValue $ws.Range("D48") '1.932'
$ws.Range("E48").Value = '  -1.05%  '

Set-TextValue $ws.Range("D49") '1.127'
$ws.Range("E49").Value = '  -0.72%  '

$ws.Range("E50").Value = '  -1.60%  '

Set-TextValue $ws.Range("D51") '72.10'
$ws.Range("E51").Value = '  -0.27%  '
